$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$ws.Range("C9").Value = 'No nationwide lockdown. 23 of Brazil''s 27 federative units (states) have imposed stronger confinement measures. Authorizations for entry into indigenous lands are suspended, with the exception of those related to essential services.'

$ws.Range("H9").Value = 'A state of emergency declaration by Congress exempts the government from meeting 2020 fiscal target. Brazil''s government passed a large emergency package, but this also includes anticipated disbursements of benefits and postponed tax liabilities. In the areas of social protection and health, however, significant additional spending measures are being taken.'

$ws.Range("I38").Value = '- Temporary compensation of wage costs, up to 90% of wage bill, conditional on at least 20% fall in turnover in the months March to May compared to 2019 (EUR 10bn). 
- In addiiton to the wage scheme, the government will cover 30% of pension contributions and the employers'' premium (1 April)
- Self-employed and independent entrepreneurs without staff can receive social assistance benefits for theree months through an accelerated procedure, where asset tests and tests on partner income are skipped. Loans at reduced interest rates are available under the same scheme (EUR 2bn).
- Self-employed (without employees) hit by the corona crisis are to receive monthly 1,050 ''euros over the next three months, irrespective of their savings. For married couples or couples with children a maximum of EUR 1,500 a month will apply (30 March). '

$ws.Range("J38").Value = '- Tax collection of VAT, income, corporate, wage and turnover taxes will be halted immediately upon application to defer tax payments from companies. The estimated size of this support is between EUR 45-60 bn (5.8%-8.4%) of GDP in 2020. Companies get more time to provide documentation. The interest rate on tax debts is reduced from 4% to almost 0%. 
- A guarantee facility for SME loans (BMKB) is already in place and working – EUR 665 million of an available budget of EUR 765 million, 0.08% of GDP. The Ministry of Economic Affairs and Climate Policy has expanded guarantees to a higher guarantee percentage, up to 90% of the credit risk of banks for SMEs that want to take out a loan of EUR 1.5 mln maximum but that do not have the required collateral. 
- The Guarantee for Entrepreneurial Financing (GO-regeling) for businesses that experience problems in obtaining bank loans and bank guarantees has been broadened. The guarantee ceiling increases from 400 million to 1.5 billion euros. Both SMEs and large companies can obtain a 50% guarantee on bank loans and bank guarantees (minimum 1.5 million - maximum 50 million euros per company). The maximum per company is temporarily increased to 150 million euros. 
-  Firms that needed to close due to the coronavurus can receive a one-off lump sum allowance of EUR 4000 (30 March). In a first instance, this included firms in the hospitality sector, hairdressers and the travel industry. On 31 March, the emergency desk was extended to retailers.
- The banking sector has responded by giving an aoutomatic 6-month payment holiday (interest and amortisation) for all business loans of less than 2.5 million euros. On 30 March, the bank ABN Amro raised the ceiling to EUR 50 million to extend the loan repayment pause to bigger enterprises. 
Interest rates are also lowered for a publicly backed microcredit scheme (Qredits).
- A temporary guarantee for working capital in agricultural and horticultural companies will be provided under the Guarantee for SME Agricultural Credits (BL) scheme.
- The government increased export credit insurance for firms through the credit insurer Atradius. The share of the working capital that companies need for their export production was raised from 80 to 95 percent (30 March)
- Competition regulations and driving and rest times regulations are eased temporarily to secure supply of foodstuffs, energy pharmaceutical products, as well as transport of waste and garbage.
'

$ws.Range("I53").Value = '• Minimum monthly pensions were raised from 1.000 to 1.500 Turkish lira (230 US$).                                                                 • The holiday bonus for retirees will be paid in April instead of May.                                                                                            • A complementary budget of 2.1 billion Turkish lira (308 million US$) is earmarked for families in need, 2.1 million families will receive 1000 Turkish Lira (154 US$).                                                                • The conditions of utilisation of the Short-Time Working Scheme have been eased:                                                                                             - The requirement of 600 days of contribution is reduced to 450 days;
- The requirement of an employment contract for the last 120 days is reduced to 60 days;
- The application procedure was simplified with lighter documentation;
- A new criterion of “not laying off workers” was added for employers willing to benefit from the Programme.
• Compensatory working time was increased from 2 months to 4 months.                                                                                     • Civil servants are encouraged to flexible and remote working, provided that the continuity of public service is ensured.                • Flexible working and remote-working  are also encouraged in the private sector. 
• Final filing and payment dates of personal income tax returns are extended.
• Taxpayers directly affected by the coronavirus outbreak will benefit from the force majeure provisions of the Tax Procedure Law.  
• Due to the curfew restrictions for citizens above 65 years and  those with chronic illnesses, the curfew period will be considered as force majeure for these taxpayers.                                          • Municipalities, including in Ankara and Istanbul, introduced solidarity schemes funded by local residents to provide cash and food aid to local households in need.                                            • On 30 March, President Erdogan launched a national solidarity campaign to help citizens in difficulty, inviting businesses and households to contribute generously. As of 6 April, 1.5 billion TL were collected in this National Sloidarity Fund.                                                       • On 7 April, additional social transfers (of 1000 TL, 154 US$) to 2.3 million households which had not qualified for the 18 March measures (with a fiscal cost of 2.3 billion TL), and to workers losing their job without qualifying for short-term working arrangements or unemployment insurance were announced.                                                                        • From 7 April, households in need ineligible for other   programmes can apply for special support from the National Solidarity Fund. '

$ws.Range("B5:B56").Value = 43930

"done"